$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22; existing rows 22-29 shift down to 23-30
# (row 29's data ends up duplicated into the new row 30, unchanged).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(22, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(22, 4).Value = 44726
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = 100112044
$ws.Cells.Item(22, 7).Value = 'Perejil'
$ws.Cells.Item(22, 8).Value = 'Sin especificar'
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 250
$ws.Cells.Item(22, 11).Value = 2500
$ws.Cells.Item(22, 12).Value = 2800
$ws.Cells.Item(22, 13).Value = 2650
$ws.Cells.Item(22, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(22, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(22, 16).Value = 1325
$ws.Cells.Item(22, 17).Value = 2
$ws.Cells.Item(22, 18).Value = 'Hortaliza'
